$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 1: "This is a Microsoft word document." gets two trailing
#    spaces, followed by three red (C00000) runs that spell out:
#    "(This is a change – Ve" + "rsion for branch alternate" + ")"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$insPoint = $d.Range($p1.End - 1, $p1.End - 1)
$insPoint.InsertAfter("  ")

$p1b = $d.Paragraphs(1).Range
$insPoint2 = $d.Range($p1b.End - 1, $p1b.End - 1)
$insPoint2.InsertAfter([char]0x0028 + "This is a change " + [char]0x2013 + " Ve")
$run1 = $d.Range($insPoint2.Start, $d.Paragraphs(1).Range.End - 1)
$run1.Font.Color = 192

$p1c = $d.Paragraphs(1).Range
$insPoint3 = $d.Range($p1c.End - 1, $p1c.End - 1)
$insPoint3.InsertAfter("rsion for branch alternate")
$run2 = $d.Range($insPoint3.Start, $d.Paragraphs(1).Range.End - 1)
$run2.Font.Color = 192

$p1d = $d.Paragraphs(1).Range
$insPoint4 = $d.Range($p1d.End - 1, $p1d.End - 1)
$insPoint4.InsertAfter([char]0x0029)
$run3 = $d.Range($insPoint4.Start, $d.Paragraphs(1).Range.End - 1)
$run3.Font.Color = 192

# ---------------------------------------------------------------------------
# 2) The empty paragraph (currently <w:p/>) right before "The Raven" gets
#    shading + paragraph-mark formatting (Calibri/bold/color), but stays
#    textually empty.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)

# Apply the properties that (per this host's Font plumbing) only stick to
# the paragraph-mark rPr when set on a *collapsed* range.
$rEmpty = $p3.Range
$rEmpty.Font.NameFarEast = "Times New Roman"
$rEmpty.Font.NameBi = "Calibri"
$rEmpty.Font.BoldBi = $true

# Apply the properties that (on this host) only stick to the paragraph-mark
# rPr when the range actually contains a character - so type a placeholder,
# format it, then remove it again.
$r = $d.Paragraphs(3).Range
$r.InsertAfter("X")
$r2 = $d.Paragraphs(3).Range
$r2.Font.Name = "Calibri"
$r2.Font.Bold = $true
$r2.Font.Color = 2236704

# Paragraph shading: clear pattern, auto foreground, F9F9F9 background.
$p3.Shading.Texture = 0
$p3.Shading.ForegroundPatternColor = -16777216
$p3.Shading.BackgroundPatternColor = 16382457

# Remove the placeholder character, leaving an empty paragraph whose pPr
# now carries the formatting applied above.
$delStart = $d.Paragraphs(3).Range.Start
$d.Range($delStart, $delStart + 1).Delete()

# ---------------------------------------------------------------------------
# 3) Final paragraph ("ank God almighty, we are free at last.") is cleared
#    out entirely, leaving a bare empty paragraph.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$pLast = $d.Paragraphs($n)
$lastRange = $pLast.Range
$d.Range($lastRange.Start, $lastRange.End - 1).Delete()
$d.Paragraphs($n).Style = "Normal"

Write-Output "done"
